# Applies the dated/value updates described by the diff.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-24 Thursday", "2025-07-25 Friday"),
    @("965×2=1930", "892×3=2676"),
    @("318×3=954", "984×4=3936"),
    @("880×4=3520", "638×2=1276"),
    @("613×4=2452", "949×3=2847"),
    @("480×8=3840", "427×9=3843"),
    @("936×6=5616", "835×2=1670"),
    @("568×9=5112", "228×9=2052"),
    @("606×9=5454", "333×9=2997"),
    @("933×5=4665", "679×6=4074"),
    @("262×5=1310", "590×3=1770"),
    @("306×6=1836", "902×9=8118"),
    @("679×9=6111", "967×5=4835"),
    @("847×9=7623", "603×9=5427"),
    @("181×2=362", "160×4=640"),
    @("783×3=2349", "114×8=912"),
    @("749×2=1498", "719×8=5752"),
    @("140×8=1120", "749×9=6741"),
    @("702×9=6318", "337×7=2359"),
    @("266×6=1596", "555×3=1665"),
    @("523×6=3138", "151×5=755"),
    @("865×8=6920", "298×9=2682"),
    @("174×9=1566", "990×6=5940"),
    @("376×2=752", "583×2=1166"),
    @("740×2=1480", "925×7=6475"),
    @("923×3=2769", "286×3=858")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$d.Save()
